# Auto-generated Excel COM-interop script
# Updates literal market-price values in the Anima Profits workbook
# (per-sheet leve profit calculations) to match refreshed Universalis data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2761.5386
$ws.Range("I40").Value = 1487.5
$ws.Range("J40").Value = 4800
$ws.Range("K40").Value = 1487.5
$ws.Range("L40").Value = 4800
$ws.Range("M40").Value = -1312.5
$ws.Range("N40").Value = -5150
$ws.Range("H86").Value = 69582170
$ws.Range("I86").Value = 80018830
$ws.Range("J86").Value = 4400
$ws.Range("K86").Value = 80018830
$ws.Range("L86").Value = 4400
$ws.Range("M86").Value = -80017707
$ws.Range("N86").Value = -6646
$ws.Range("H89").Value = 69582170
$ws.Range("I89").Value = 80018830
$ws.Range("J89").Value = 4400
$ws.Range("K89").Value = 400094150
$ws.Range("L89").Value = 22000
$ws.Range("M89").Value = -400088534
$ws.Range("N89").Value = -33232
$ws.Range("H92").Value = 8548103
$ws.Range("I92").Value = 13333850
$ws.Range("J92").Value = 2127.7856
$ws.Range("K92").Value = 13333850
$ws.Range("L92").Value = 2127.7856
$ws.Range("M92").Value = -13332602
$ws.Range("N92").Value = -4623.7856
$ws.Range("H106").Value = 20692344
$ws.Range("I106").Value = 27275476
$ws.Range("K106").Value = 27275476
$ws.Range("M106").Value = -27274845
$ws.Range("H137").Value = 1683.48
$ws.Range("I137").Value = 1331.2307
$ws.Range("K137").Value = 3993.6921
$ws.Range("M137").Value = -1443.6921
$ws.Range("H141").Value = 6164.2
$ws.Range("I141").Value = 2931
$ws.Range("J141").Value = 11014
$ws.Range("K141").Value = 8793
$ws.Range("L141").Value = 33042
$ws.Range("M141").Value = -3613
$ws.Range("N141").Value = -43402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 640.1667
$ws.Range("I2").Value = 695.7143
$ws.Range("J2").Value = 562.4
$ws.Range("K2").Value = 695.7143
$ws.Range("L2").Value = 562.4
$ws.Range("M2").Value = -582.7143
$ws.Range("N2").Value = -788.4
$ws.Range("H45").Value = 2443.7693
$ws.Range("I45").Value = 2087.5
$ws.Range("J45").Value = 3013.8
$ws.Range("K45").Value = 2087.5
$ws.Range("L45").Value = 3013.8
$ws.Range("M45").Value = -1710.5
$ws.Range("N45").Value = -3767.8
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16450
$ws.Range("H116").Value = 640.1667
$ws.Range("I116").Value = 695.7143
$ws.Range("J116").Value = 562.4
$ws.Range("K116").Value = 695.7143
$ws.Range("L116").Value = 562.4
$ws.Range("M116").Value = 1598.2857
$ws.Range("N116").Value = -5150.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 640.1667
$ws.Range("I3").Value = 695.7143
$ws.Range("J3").Value = 562.4
$ws.Range("K3").Value = 695.7143
$ws.Range("L3").Value = 562.4
$ws.Range("M3").Value = -581.7143
$ws.Range("N3").Value = -790.4
$ws.Range("H134").Value = 2876.721
$ws.Range("I134").Value = 2877.4194
$ws.Range("J134").Value = 2874.9167
$ws.Range("K134").Value = 8632.2582
$ws.Range("L134").Value = 8624.750100000001
$ws.Range("M134").Value = -6097.2582
$ws.Range("N134").Value = -13694.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2092.2856
$ws.Range("I99").Value = 1498.4
$ws.Range("K99").Value = 1498.4
$ws.Range("M99").Value = -0.4000000000000909
$ws.Range("H107").Value = 3677328.8
$ws.Range("I107").Value = 6250511
$ws.Range("K107").Value = 6250511
$ws.Range("M107").Value = -6248591
$ws.Range("H126").Value = 2092.2856
$ws.Range("I126").Value = 1498.4
$ws.Range("K126").Value = 4495.200000000001
$ws.Range("M126").Value = -2025.200000000001
$ws.Range("H134").Value = 15225.5
$ws.Range("I134").Value = 19317.334
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 57952.00199999999
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -55417.00199999999
$ws.Range("N134").Value = -13920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2474.3965
$ws.Range("I122").Value = 445.5909
$ws.Range("K122").Value = 4010.3181
$ws.Range("M122").Value = -1560.3181
$ws.Range("H138").Value = 7985.7144
$ws.Range("I138").Value = 3200
$ws.Range("J138").Value = 19950
$ws.Range("K138").Value = 9600
$ws.Range("L138").Value = 59850
$ws.Range("M138").Value = -4460
$ws.Range("N138").Value = -70130
$ws.Range("H140").Value = 1509.1052
$ws.Range("I140").Value = 1341.1111
$ws.Range("J140").Value = 4533
$ws.Range("K140").Value = 4023.3333
$ws.Range("L140").Value = 13599
$ws.Range("M140").Value = 1156.6667
$ws.Range("N140").Value = -23959
$ws.Range("H141").Value = 9834.388999999999
$ws.Range("I141").Value = 8913.223
$ws.Range("J141").Value = 10755.556
$ws.Range("K141").Value = 26739.669
$ws.Range("L141").Value = 32266.668
$ws.Range("M141").Value = -21559.669
$ws.Range("N141").Value = -42626.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H132").Value = 2481.366
$ws.Range("I132").Value = 1982.2084
$ws.Range("J132").Value = 3186.0588
$ws.Range("K132").Value = 5946.6252
$ws.Range("L132").Value = 9558.1764
$ws.Range("M132").Value = -3416.6252
$ws.Range("N132").Value = -14618.1764
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12650.066
$ws.Range("I22").Value = 635.5714
$ws.Range("J22").Value = 23162.75
$ws.Range("K22").Value = 635.5714
$ws.Range("L22").Value = 23162.75
$ws.Range("M22").Value = -340.5714
$ws.Range("N22").Value = -23752.75
$ws.Range("H27").Value = 12650.066
$ws.Range("I27").Value = 635.5714
$ws.Range("J27").Value = 23162.75
$ws.Range("K27").Value = 635.5714
$ws.Range("L27").Value = 23162.75
$ws.Range("M27").Value = -528.5714
$ws.Range("N27").Value = -23376.75
$ws.Range("H93").Value = 11522.091
$ws.Range("I93").Value = 13249.223
$ws.Range("J93").Value = 3750
$ws.Range("K93").Value = 13249.223
$ws.Range("L93").Value = 3750
$ws.Range("M93").Value = -12001.223
$ws.Range("N93").Value = -6246
$ws.Range("H122").Value = 5497.375
$ws.Range("I122").Value = 4489.5
$ws.Range("K122").Value = 13468.5
$ws.Range("M122").Value = -11018.5
